# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the Membrillo price rows: the values in
# columns D, L, M, N, O, P, Q, R, S, T for each data row (2..30) are
# replaced by the values that previously lived in a different row of the
# same block (columns A, B, C, E, F, G, H, I, J, K never change). Capture
# every source row first, then write the permuted values back so a cell
# is never read after it has already been overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new-row -> old-row the replacement values are sourced from
$perm = @{
    2  = 26
    3  = 27
    4  = 23
    5  = 7
    6  = 8
    7  = 21
    8  = 22
    9  = 18
    10 = 19
    11 = 24
    12 = 25
    13 = 11
    14 = 12
    15 = 13
    16 = 14
    17 = 20
    18 = 28
    19 = 29
    20 = 9
    21 = 10
    22 = 3
    23 = 4
    24 = 15
    25 = 16
    26 = 17
    27 = 5
    28 = 6
    29 = 30
    30 = 2
}

# columns (1-based) that participate in the shuffle
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D L M N O P Q R S T

# snapshot every source cell before any writes happen (Value2 = raw
# underlying value, unlike Value which can come back formatted)
$snapshot = @{}
foreach ($r in 2..30) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# write the permuted values back
foreach ($r in 2..30) {
    $src = $perm[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
